# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the bf4a510b-5d1c-4f1e-96ae-c1ad50c1ab36.md row to reflect
# the latest handoff XLIFF generation.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# bf4a510b-... row (row 7).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-05 20:48:48"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the bf4a510b-... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-05 20:48:43"

# de-de sheet: "Latest Handoff Datetime" column (H) for the bf4a510b-... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-05 20:48:48"
